$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("B6").Value = 6750018
$ws.Range("E6").Value = "WDA Swiecie"
$ws.Range("F6").Value = "Swit Starozreby"
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("I6").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("K6").Value = "D"
$ws.Range("L6").Value = 1.571
$ws.Range("M6").Value = 4
$ws.Range("N6").Value = 4.333
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 4.5
$ws.Range("Q6").Value = 6
$ws.Range("R6").Value = -1.25
$ws.Range("S6").Value = 1.8
$ws.Range("T6").Value = 2
$ws.Range("U6").Value = 3.25
$ws.Range("V6").Value = 1.8
$ws.Range("W6").Value = 2
$ws.Range("X6").Value = -1
$ws.Range("Y6").Value = 3.5
$ws.Range("Z6").Value = -1
$ws.Range("AA6").Value = -1
$ws.Range("AB6").Value = 1
$ws.Range("AC6").Value = -1
$ws.Range("AD6").Value = 1

# Row 7
$ws.Range("B7").Value = 6746871
$ws.Range("E7").Value = "Hutnik Warsaw"
$ws.Range("F7").Value = "Chemik Bydgoszcz"
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 2
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = "H"
$ws.Range("L7").Value = 3.25
$ws.Range("M7").Value = 4
$ws.Range("N7").Value = 1.8
$ws.Range("O7").Value = 3.5
$ws.Range("P7").Value = 4
$ws.Range("Q7").Value = 1.727
$ws.Range("R7").Value = 0.75
$ws.Range("S7").Value = 1.825
$ws.Range("T7").Value = 1.975
$ws.Range("U7").Value = 3.25
$ws.Range("V7").Value = 1.9
$ws.Range("W7").Value = 1.9
$ws.Range("X7").Value = 2.5
$ws.Range("Y7").Value = -1
$ws.Range("Z7").Value = -1
$ws.Range("AA7").Value = 0.825
$ws.Range("AB7").Value = -1
$ws.Range("AC7").Value = -1
$ws.Range("AD7").Value = 0.8999999999999999

# Row 11
$ws.Range("B11").Value = 6757267
$ws.Range("E11").Value = "Czarni Pruszcz Gdanski"
$ws.Range("F11").Value = "Chojniczanka Chojnice II"
$ws.Range("G11").Value = 4
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = "H"
$ws.Range("L11").Value = 1.909
$ws.Range("M11").Value = 4
$ws.Range("N11").Value = 3
$ws.Range("O11").Value = 1.909
$ws.Range("P11").Value = 4
$ws.Range("Q11").Value = 3
$ws.Range("R11").Value = -0.5
$ws.Range("S11").Value = 1.95
$ws.Range("T11").Value = 1.85
$ws.Range("U11").Value = 3.75
$ws.Range("V11").Value = 2
$ws.Range("W11").Value = 1.8
$ws.Range("X11").Value = 0.909
$ws.Range("Y11").Value = -1
$ws.Range("Z11").Value = -1
$ws.Range("AA11").Value = 0.95
$ws.Range("AB11").Value = -1
$ws.Range("AC11").Value = 0.5
$ws.Range("AD11").Value = -0.5

# Row 12
$ws.Range("B12").Value = 6761111
$ws.Range("E12").Value = "Arka Gdynia II"
$ws.Range("F12").Value = "Sparta Sycewice"
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = "D"
$ws.Range("L12").Value = 2
$ws.Range("M12").Value = 4.333
$ws.Range("N12").Value = 2.625
$ws.Range("O12").Value = 2
$ws.Range("P12").Value = 4.333
$ws.Range("Q12").Value = 2.625
$ws.Range("R12").Value = -0.25
$ws.Range("S12").Value = 1.85
$ws.Range("T12").Value = 1.95
$ws.Range("U12").Value = 3.25
$ws.Range("V12").Value = 1.85
$ws.Range("W12").Value = 1.95
$ws.Range("X12").Value = -1
$ws.Range("Y12").Value = 3.333
$ws.Range("Z12").Value = -1
$ws.Range("AA12").Value = -0.5
$ws.Range("AB12").Value = 0.475
$ws.Range("AC12").Value = -1
$ws.Range("AD12").Value = 0.95

# Row 183
$ws.Range("B183").Value = 8229582
$ws.Range("E183").Value = "Wisla Krakow II"
$ws.Range("F183").Value = "Dalin Myslenice"
$ws.Range("G183").Value = 2
$ws.Range("H183").Value = 0
$ws.Range("I183").Value = 0
$ws.Range("J183").Value = 0
$ws.Range("K183").Value = "H"
$ws.Range("L183").Value = 1.2
$ws.Range("M183").Value = 6.5
$ws.Range("N183").Value = 8
$ws.Range("O183").Value = 1.285
$ws.Range("P183").Value = 5.5
$ws.Range("Q183").Value = 6.25
$ws.Range("R183").Value = -1.75
$ws.Range("S183").Value = 2
$ws.Range("T183").Value = 1.8
$ws.Range("U183").Value = 3
$ws.Range("V183").Value = 1.775
$ws.Range("W183").Value = 2.025
$ws.Range("X183").Value = 0.2849999999999999
$ws.Range("Y183").Value = -1
$ws.Range("Z183").Value = -1
$ws.Range("AA183").Value = 0.5
$ws.Range("AB183").Value = -0.5
$ws.Range("AC183").Value = -1
$ws.Range("AD183").Value = 1.025

# Row 184
$ws.Range("B184").Value = 8229175
$ws.Range("E184").Value = "Wierna Malogoszcz"
$ws.Range("F184").Value = "Klimontowianka Klimontow"
$ws.Range("G184").Value = 2
$ws.Range("H184").Value = 3
$ws.Range("I184").Value = 1
$ws.Range("J184").Value = 0
$ws.Range("K184").Value = "A"
$ws.Range("L184").Value = 2.25
$ws.Range("M184").Value = 3.75
$ws.Range("N184").Value = 2.5
$ws.Range("O184").Value = 2.25
$ws.Range("P184").Value = 3.75
$ws.Range("Q184").Value = 2.5
$ws.Range("R184").Value = 0
$ws.Range("S184").Value = 1.8
$ws.Range("T184").Value = 2
$ws.Range("U184").Value = 3
$ws.Range("V184").Value = 1.8
$ws.Range("W184").Value = 2
$ws.Range("X184").Value = -1
$ws.Range("Y184").Value = -1
$ws.Range("Z184").Value = 1.5
$ws.Range("AA184").Value = -1
$ws.Range("AB184").Value = 1
$ws.Range("AC184").Value = 0.8
$ws.Range("AD184").Value = -1

# Row 185
$ws.Range("B185").Value = 8229174
$ws.Range("E185").Value = "Sparta Katowice"
$ws.Range("F185").Value = "Znicz Klobuck"
$ws.Range("G185").Value = 2
$ws.Range("H185").Value = 0
$ws.Range("I185").Value = 1
$ws.Range("J185").Value = 0
$ws.Range("K185").Value = "H"
$ws.Range("L185").Value = 1.25
$ws.Range("M185").Value = 5.5
$ws.Range("N185").Value = 9
$ws.Range("O185").Value = 1.25
$ws.Range("P185").Value = 5.5
$ws.Range("Q185").Value = 9
$ws.Range("R185").Value = -1.75
$ws.Range("S185").Value = 1.8
$ws.Range("T185").Value = 2
$ws.Range("U185").Value = 3.75
$ws.Range("V185").Value = 1.875
$ws.Range("W185").Value = 1.925
$ws.Range("X185").Value = 0.25
$ws.Range("Y185").Value = -1
$ws.Range("Z185").Value = -1
$ws.Range("AA185").Value = 0.4
$ws.Range("AB185").Value = -0.5
$ws.Range("AC185").Value = -1
$ws.Range("AD185").Value = 0.925

# Row 222
$ws.Range("B222").Value = 8284778
$ws.Range("E222").Value = "DKS Dobre Miasto"
$ws.Range("F222").Value = "Mragowia Mragowo"
$ws.Range("G222").Value = 2
$ws.Range("H222").Value = 0
$ws.Range("I222").ClearContents()
$ws.Range("J222").ClearContents()
$ws.Range("K222").Value = "H"
$ws.Range("L222").Value = 2.45
$ws.Range("M222").Value = 3.7
$ws.Range("N222").Value = 2.3
$ws.Range("O222").Value = 2.9
$ws.Range("P222").Value = 3.75
$ws.Range("Q222").Value = 2
$ws.Range("R222").Value = 0.25
$ws.Range("S222").Value = 1.975
$ws.Range("T222").Value = 1.825
$ws.Range("U222").Value = 3.25
$ws.Range("V222").Value = 1.9
$ws.Range("W222").Value = 1.9
$ws.Range("X222").Value = 1.9
$ws.Range("Y222").Value = -1
$ws.Range("Z222").Value = -1
$ws.Range("AA222").Value = 0.9750000000000001
$ws.Range("AB222").Value = -1
$ws.Range("AC222").Value = -1
$ws.Range("AD222").Value = 0.8999999999999999

# Row 223
$ws.Range("B223").Value = 8287558
$ws.Range("E223").Value = "WDA Swiecie"
$ws.Range("F223").Value = "Mazovia Minsk Mazowiecki"
$ws.Range("G223").Value = 1
$ws.Range("H223").Value = 2
$ws.Range("I223").ClearContents()
$ws.Range("J223").ClearContents()
$ws.Range("K223").Value = "A"
$ws.Range("L223").Value = 1.727
$ws.Range("M223").Value = 3.5
$ws.Range("N223").Value = 4
$ws.Range("O223").Value = 1.727
$ws.Range("P223").Value = 3.5
$ws.Range("Q223").Value = 4
$ws.Range("R223").Value = -0.75
$ws.Range("S223").Value = 1.95
$ws.Range("T223").Value = 1.85
$ws.Range("U223").Value = 3.5
$ws.Range("V223").Value = 1.825
$ws.Range("W223").Value = 1.975
$ws.Range("X223").Value = -1
$ws.Range("Y223").Value = -1
$ws.Range("Z223").Value = 3
$ws.Range("AA223").Value = -1
$ws.Range("AB223").Value = 0.8500000000000001
$ws.Range("AC223").Value = -1
$ws.Range("AD223").Value = 0.9750000000000001

